$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update individual odds values on existing rows ---
$ws.Range("L2").Value = 1.41
$ws.Range("AB2").Value = 13
$ws.Range("AD2").Value = 12.5

$ws.Range("AH3").Value = 22

$ws.Range("F4").Value = 2.38
$ws.Range("P4").Value = 2.08
$ws.Range("S4").Value = 2.84
$ws.Range("T4").Value = 1.63
$ws.Range("U4").Value = 2.28

$ws.Range("F5").Value = 8.6
$ws.Range("G5").Value = 8.800000000000001
$ws.Range("H5").Value = 1.44
$ws.Range("I5").Value = 1.45
$ws.Range("K5").Value = 5.3
$ws.Range("N5").Value = 5
$ws.Range("P5").Value = 2.38
$ws.Range("Q5").Value = 1.7
$ws.Range("R5").Value = 1.53
$ws.Range("S5").Value = 2.8
$ws.Range("T5").Value = 1.93
$ws.Range("V5").Value = 3.2
$ws.Range("W5").Value = 1.12
$ws.Range("X5").Value = 22
$ws.Range("Y5").Value = 9.199999999999999
$ws.Range("AJ5").Value = 280
$ws.Range("AO5").Value = 6.2

$ws.Range("L6").Value = 1.28

$ws.Range("Q7").Value = 1.76
$ws.Range("AK7").Value = 48

$ws.Range("F9").Value = 2.34

$ws.Range("G11").Value = 2.74
$ws.Range("P11").Value = 2
$ws.Range("R11").Value = 1.4
$ws.Range("W11").Value = 1.57

$ws.Range("I12").Value = 1.97
$ws.Range("R12").Value = 1.54

$ws.Range("G13").Value = 5
$ws.Range("I13").Value = 1.8
$ws.Range("T13").Value = 1.6

# --- Append a new row (14) with a new match ---
$ws.Range("A14").Value = "Colombian Primera A"
# Force the date-like text to remain plain text (avoid Excel auto-converting
# it to a date serial number) by using a leading quote prefix.
$ws.Range("B14").Value = "'2026-02-18"
$ws.Range("C14").Value = "21:30:00"
$ws.Range("D14").Value = "Junior FC Barranquilla"
$ws.Range("E14").Value = "America de Cali S.A"
$ws.Range("F14").Value = 2.04
$ws.Range("G14").Value = 2.1
$ws.Range("H14").Value = 4
$ws.Range("I14").Value = 4.5
$ws.Range("J14").Value = 3.4
$ws.Range("K14").Value = 3.6
$ws.Range("L14").Value = 1.01
$ws.Range("M14").Value = 1.09
$ws.Range("N14").Value = 2.52
$ws.Range("O14").Value = 1.41
$ws.Range("P14").Value = 1.7
$ws.Range("Q14").Value = 2.26
$ws.Range("R14").Value = 1.2
$ws.Range("S14").Value = 3.7
$ws.Range("T14").Value = 1.96
$ws.Range("U14").Value = 1.78
$ws.Range("V14").Value = 1.27
$ws.Range("W14").Value = 1.9
$ws.Range("X14").Value = 1000
$ws.Range("Y14").Value = 14.5
$ws.Range("Z14").Value = 1000
$ws.Range("AA14").Value = 120
$ws.Range("AB14").Value = 8.800000000000001
$ws.Range("AC14").Value = 8.800000000000001
$ws.Range("AD14").Value = 21
$ws.Range("AE14").Value = 1000
$ws.Range("AF14").Value = 13.5
$ws.Range("AG14").Value = 12.5
$ws.Range("AH14").Value = 1000
$ws.Range("AI14").Value = 1000
$ws.Range("AJ14").Value = 1000
$ws.Range("AK14").Value = 1000
$ws.Range("AL14").Value = 1000
$ws.Range("AM14").Value = 170
$ws.Range("AN14").Value = 24
$ws.Range("AO14").Value = 1000
